$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column headers (I0, IF) copying the bold/bordered header style from H1
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value2 = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value2 = "IF"

# Fill I and J columns for rows 2-23: I is always 1, J mirrors column H
for ($r = 2; $r -le 23; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
